# Applies the LOB1039 worksheet restructuring described in the commit diff:
#  - sharedStrings gain 3 new entries (Portuguese Objectives text, PT summary,
#    PT full syllabus) and reorder others
#  - worksheet rows 13-23 are rebuilt/shifted by one, row 24 is newly added
#  - column A's <col> span is narrowed from A:B to just A (fixed as a side
#    effect of writing column A explicitly in every row below)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Objetivos:) - text-only correction, no structural change
$ws.Cells.Item(10, 2).Value = 'Observação experimental de fenômenos relacionados à eletricidade e magnetismo.'
$ws.Cells.Item(10, 3).Value = 'Observação experimental de fenômenos relacionados à eletricidade e magnetismo.'

# Rows 13-23 are fully rebuilt (content/styles/heights all change); clear each
# one first so no stale cell/style survives from the old layout
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(23).Insert()
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(13).Insert()

# Row 13
$ws.Cells.Item(13, 2).Value = '3268262 - Carlos Renato Menegatti'
$ws.Cells.Item(13, 3).Value = '3268262 - Carlos Renato Menegatti'
$ws.Cells.Item(13, 1).Clear()
$ws.Cells.Item(10, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial(-4122)

# Row 14
$ws.Cells.Item(14, 1).Value = 'Programa resumido:'
$ws.Cells.Item(14, 2).Value = 'Campo Eletrostático e Mapeamento de Equipotenciais;  Introdução a Circuitos de Corrente Contínua; Resistência, Resistividade e Corrente Elétrica; Circuitos de Corrente Contínua;  Capacitores; Voltímetros, Amperímetros e Ohmímetros; Osciloscópios; Campo Magnetostático; Lei de Indução de Faraday; Circuitos RL e RC;'
$ws.Cells.Item(14, 3).Value = 'Campo Eletrostático e Mapeamento de Equipotenciais;  Introdução a Circuitos de Corrente Contínua; Resistência, Resistividade e Corrente Elétrica; Circuitos de Corrente Contínua;  Capacitores; Voltímetros, Amperímetros e Ohmímetros; Osciloscópios; Campo Magnetostático; Lei de Indução de Faraday; Circuitos RL e RC;'
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Cells.Item(15, 1).Value = 'Short syllabus:'
$ws.Cells.Item(15, 2).Value = 'Field Electrostatic Equipotential Mapping; Introduction to Direct Current Circuits; Resistance, Resistivity and Electrical Current; Kirchoff Laws; Capacitors; Voltmeters, Ammeters and Ohmmeters; Oscilloscopes; Magnetostatic Field; Faraday''s Induction Law; Circuits RL and RC;'
$ws.Cells.Item(15, 3).Value = 'Field Electrostatic Equipotential Mapping; Introduction to Direct Current Circuits; Resistance, Resistivity and Electrical Current; Kirchoff Laws; Capacitors; Voltmeters, Ammeters and Ohmmeters; Oscilloscopes; Magnetostatic Field; Faraday''s Induction Law; Circuits RL and RC;'
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Cells.Item(16, 1).Value = 'Programa:'
$ws.Cells.Item(16, 2).Value = '1) Campo Eletrostático e Mapeamento de Equipotenciais: Campo de placas paralelas, Campo de cargas pontuais, Efeito de isolante e condutor.2) Introdução a Circuitos de Corrente Contínua: Resistores ôhmicos, Resistores não-ohmicos.3) Resistência e Corrente Elétrica: Lei de Ohm, Modelo de Drude.4) Circuitos de Corrente Contínua: Leis de Kirchoff.5) Capacitores: Associação de capacitores, Carga e descarga de um capacitor.6) Voltímetros, Amperímetros e Ohmímetros: Princípio de funcionamento do Galvanômetro, Construção de Voltímetros, Amperímetros e Ohmímetros.7) Osciloscópios: Princípio de Funcionamento do Osciloscópio.8) Campo Magnetostático: Lei de Biot-Savart, Lei de Ampère, Efeito Hall.9) Lei de Indução de Faraday: Indutância mútua e auto-indutância, Geração de tensão AC.10) Circuitos RL e RC em corrente contínua.'
$ws.Cells.Item(16, 3).Value = '1) Campo Eletrostático e Mapeamento de Equipotenciais: Campo de placas paralelas, Campo de cargas pontuais, Efeito de isolante e condutor.2) Introdução a Circuitos de Corrente Contínua: Resistores ôhmicos, Resistores não-ohmicos.3) Resistência e Corrente Elétrica: Lei de Ohm, Modelo de Drude.4) Circuitos de Corrente Contínua: Leis de Kirchoff.5) Capacitores: Associação de capacitores, Carga e descarga de um capacitor.6) Voltímetros, Amperímetros e Ohmímetros: Princípio de funcionamento do Galvanômetro, Construção de Voltímetros, Amperímetros e Ohmímetros.7) Osciloscópios: Princípio de Funcionamento do Osciloscópio.8) Campo Magnetostático: Lei de Biot-Savart, Lei de Ampère, Efeito Hall.9) Lei de Indução de Faraday: Indutância mútua e auto-indutância, Geração de tensão AC.10) Circuitos RL e RC em corrente contínua.'
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Cells.Item(17, 1).Value = 'Syllabus:'
$ws.Cells.Item(17, 2).Value = '1) Electrostatic Field and Equipotential Mapping: Parallel plates Field, A point charge Field, insulating effect and conductor.2) Ohm’s Law: ohmic resistors, resistors non-ohmic.3) Resistance and Electric current: Ohm''s Law, Drude model.4) Direct Current Circuits: Kirchoff laws.5) Capacitors: Capacitors association, load and discharge a capacitor.6) Voltmeters, Ammeters and ohmmeters: Galvanometer operation principle, Voltmeters Construction, Ammeters and ohmmeters.7) Oscilloscope: Oscilloscope Operation Principle.8) Magnetostatic Field: Biot-Savart law, Ampere''s law, Hall effect.9) Faraday''s Law of Induction: Mutual inductance and self-inductance, AC voltage generation.10) RL and RC in DC circuits'
$ws.Cells.Item(17, 3).Value = '1) Electrostatic Field and Equipotential Mapping: Parallel plates Field, A point charge Field, insulating effect and conductor.2) Ohm’s Law: ohmic resistors, resistors non-ohmic.3) Resistance and Electric current: Ohm''s Law, Drude model.4) Direct Current Circuits: Kirchoff laws.5) Capacitors: Capacitors association, load and discharge a capacitor.6) Voltmeters, Ammeters and ohmmeters: Galvanometer operation principle, Voltmeters Construction, Ammeters and ohmmeters.7) Oscilloscope: Oscilloscope Operation Principle.8) Magnetostatic Field: Biot-Savart law, Ampere''s law, Hall effect.9) Faraday''s Law of Induction: Mutual inductance and self-inductance, AC voltage generation.10) RL and RC in DC circuits'
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Cells.Item(18, 1).Value = 'Avaliação:'

# Row 19
$ws.Cells.Item(19, 1).Value = 'Método:'
$ws.Cells.Item(19, 2).Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Cells.Item(19, 3).Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Cells.Item(20, 1).Value = 'Critério:'
$ws.Cells.Item(20, 2).Value = 'NF≥ 5,0.'
$ws.Cells.Item(20, 3).Value = 'NF≥ 5,0.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Cells.Item(21, 1).Value = 'Norma de recuperação:'
$ws.Cells.Item(21, 2).Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Cells.Item(21, 3).Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Cells.Item(22, 1).Value = 'Bibliografia:'
$ws.Cells.Item(22, 2).Value = '1. Apostilas do Laboratório de Ensino de Física do IFSC/USP.2. VUOLO, J.H. Fundamentos da Teoria de Erros, Edgard Blucher (1996).3. NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 3, Edgard Blucher (2008).4. RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 3, LTC (2008).5. TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 3, LTC (2008).6. SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física III, Vol. 3,     Pearson Addison Wesley (2009).7. JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 3, Thomson Pioneira (2008).'
$ws.Cells.Item(22, 3).Value = '1. Apostilas do Laboratório de Ensino de Física do IFSC/USP.2. VUOLO, J.H. Fundamentos da Teoria de Erros, Edgard Blucher (1996).3. NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 3, Edgard Blucher (2008).4. RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 3, LTC (2008).5. TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 3, LTC (2008).6. SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física III, Vol. 3,     Pearson Addison Wesley (2009).7. JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 3, Thomson Pioneira (2008).'
$ws.Rows.Item(22).RowHeight = 120

# Row 23
$ws.Cells.Item(23, 1).Value = 'Requisitos:'

# Row 24 is brand new (didn't exist before)
$ws.Cells.Item(24, 2).Value = 'LOB1038 -  Física Experimental I  (Requisito fraco)
'
$ws.Cells.Item(24, 3).Value = 'LOB1038 -  Física Experimental I  (Requisito fraco)
'
$ws.Cells.Item(10, 2).Copy()
$ws.Cells.Item(24, 2).PasteSpecial(-4122)
$ws.Rows.Item(24).RowHeight = 30

$excel.CutCopyMode = $false
